$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New timelog rows (20/21/22 Apr 2015 work on "thesis chapter fundamentals") ---

# Row 95
$ws.Range("A95").Value() = 42115
$ws.Range("A95").NumberFormat() = "yyyy\-mm\-dd;@"
$ws.Range("B95").Value() = 0.47916666666666669
$ws.Range("B95").NumberFormat() = "h:mm"
$ws.Range("C95").Value() = 0.51041666666666663
$ws.Range("C95").NumberFormat() = "h:mm"
$ws.Range("D95").Formula() = "=C95-B95"
$ws.Range("D95").NumberFormat() = "hh:mm;@"
$ws.Range("E95").Value() = "thesis chapter fundamentals"

# Row 96
$ws.Range("A96").Value() = 42115
$ws.Range("A96").NumberFormat() = "yyyy\-mm\-dd;@"
$ws.Range("B96").Value() = 0.55208333333333337
$ws.Range("B96").NumberFormat() = "h:mm"
$ws.Range("C96").Value() = 0.77083333333333337
$ws.Range("C96").NumberFormat() = "h:mm"
$ws.Range("D96").Formula() = "=C96-B96"
$ws.Range("D96").NumberFormat() = "hh:mm;@"
$ws.Range("E96").Value() = "thesis chapter fundamentals"

# Row 97 (no duration formula in the source)
$ws.Range("A97").Value() = 42116
$ws.Range("A97").NumberFormat() = "yyyy\-mm\-dd;@"
$ws.Range("B97").Value() = 0.54166666666666663
$ws.Range("B97").NumberFormat() = "h:mm"
$ws.Range("C97").Value() = 0.77083333333333337
$ws.Range("C97").NumberFormat() = "h:mm"
$ws.Range("E97").Value() = "thesis chapter fundamentals"

# Rows 98-108: blank placeholder rows, only column A carries the date style
$blankRows = 98..108
foreach ($r in $blankRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat() = "yyyy\-mm\-dd;@"
}

# --- Selection bookkeeping to match the saved workbook state ---
$ws.Range("C104").Select()
